$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2275.2986
$ws.Cells.Item(15, 9).Value = 2275.2986
$ws.Cells.Item(15, 11).Value = 6825.8958
$ws.Cells.Item(15, 13).Value = -6656.8958
$ws.Cells.Item(28, 8).Value = 1097
$ws.Cells.Item(28, 9).Value = 456.57144
$ws.Cells.Item(28, 10).Value = 1737.4286
$ws.Cells.Item(28, 11).Value = 456.57144
$ws.Cells.Item(28, 12).Value = 1737.4286
$ws.Cells.Item(28, 13).Value = 28.42856
$ws.Cells.Item(28, 14).Value = -2707.4286
$ws.Cells.Item(106, 8).Value = 3690.2727
$ws.Cells.Item(106, 9).Value = 4206.615
$ws.Cells.Item(106, 10).Value = 2944.4443
$ws.Cells.Item(106, 11).Value = 4206.615
$ws.Cells.Item(106, 12).Value = 2944.4443
$ws.Cells.Item(106, 13).Value = -3575.615
$ws.Cells.Item(106, 14).Value = -4206.4443
$ws.Cells.Item(112, 8).Value = 1520.6842
$ws.Cells.Item(112, 10).Value = 1577.3889
$ws.Cells.Item(112, 12).Value = 4732.1667
$ws.Cells.Item(112, 14).Value = -6948.1667
$ws.Cells.Item(134, 8).Value = 42753.332
$ws.Cells.Item(134, 10).Value = 42753.332
$ws.Cells.Item(134, 12).Value = 42753.332
$ws.Cells.Item(134, 14).Value = -52893.332
$ws.Cells.Item(137, 8).Value = 2077.8096
$ws.Cells.Item(137, 9).Value = 1827.6444
$ws.Cells.Item(137, 11).Value = 5482.933199999999
$ws.Cells.Item(137, 13).Value = -2932.933199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 19700
$ws.Cells.Item(6, 9).Value = 19700
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 19700
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -19527
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 14).ClearContents()
$ws.Cells.Item(45, 8).Value = 1871.2941
$ws.Cells.Item(45, 9).Value = 1842.8572
$ws.Cells.Item(45, 11).Value = 1842.8572
$ws.Cells.Item(45, 13).Value = -1465.8572
$ws.Cells.Item(61, 8).Value = 5342.622
$ws.Cells.Item(61, 9).Value = 3477.3142
$ws.Cells.Item(61, 10).Value = 11871.2
$ws.Cells.Item(61, 11).Value = 3477.3142
$ws.Cells.Item(61, 12).Value = 11871.2
$ws.Cells.Item(61, 13).Value = -3265.3142
$ws.Cells.Item(61, 14).Value = -12295.2
$ws.Cells.Item(136, 8).Value = 5342.622
$ws.Cells.Item(136, 9).Value = 3477.3142
$ws.Cells.Item(136, 10).Value = 11871.2
$ws.Cells.Item(136, 11).Value = 10431.9426
$ws.Cells.Item(136, 12).Value = 35613.60000000001
$ws.Cells.Item(136, 13).Value = -7881.942599999998
$ws.Cells.Item(136, 14).Value = -40713.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 7533.3335
$ws.Cells.Item(37, 9).Value = 3000
$ws.Cells.Item(37, 11).Value = 3000
$ws.Cells.Item(37, 13).Value = -2863
$ws.Cells.Item(94, 8).Value = 1128.6875
$ws.Cells.Item(94, 9).Value = 1187.1818
$ws.Cells.Item(94, 11).Value = 1187.1818
$ws.Cells.Item(94, 13).Value = -736.1818000000001
$ws.Cells.Item(134, 8).Value = 2328.9211
$ws.Cells.Item(134, 9).Value = 2243.5925
$ws.Cells.Item(134, 10).Value = 2538.3635
$ws.Cells.Item(134, 11).Value = 6730.7775
$ws.Cells.Item(134, 12).Value = 7615.0905
$ws.Cells.Item(134, 13).Value = -4195.7775
$ws.Cells.Item(134, 14).Value = -12685.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 14650
$ws.Cells.Item(4, 9).Value = 9500
$ws.Cells.Item(4, 11).Value = 9500
$ws.Cells.Item(4, 13).Value = -9388
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 4632.643
$ws.Cells.Item(31, 9).Value = 5180.731
$ws.Cells.Item(31, 10).Value = 3742
$ws.Cells.Item(31, 11).Value = 5180.731
$ws.Cells.Item(31, 12).Value = 3742
$ws.Cells.Item(31, 13).Value = -4885.731
$ws.Cells.Item(31, 14).Value = -4332
$ws.Cells.Item(34, 8).Value = 4632.643
$ws.Cells.Item(34, 9).Value = 5180.731
$ws.Cells.Item(34, 10).Value = 3742
$ws.Cells.Item(34, 11).Value = 5180.731
$ws.Cells.Item(34, 12).Value = 3742
$ws.Cells.Item(34, 13).Value = -4978.731
$ws.Cells.Item(34, 14).Value = -4146
$ws.Cells.Item(58, 8).Value = 1569161.9
$ws.Cells.Item(58, 9).Value = 2273883
$ws.Cells.Item(58, 10).Value = 3114.5
$ws.Cells.Item(58, 11).Value = 2273883
$ws.Cells.Item(58, 12).Value = 3114.5
$ws.Cells.Item(58, 13).Value = -2273680
$ws.Cells.Item(58, 14).Value = -3520.5
$ws.Cells.Item(94, 8).Value = 1367.7142
$ws.Cells.Item(94, 9).Value = 1204
$ws.Cells.Item(94, 11).Value = 1204
$ws.Cells.Item(94, 13).Value = -753
$ws.Cells.Item(132, 8).Value = 2294.6904
$ws.Cells.Item(132, 9).Value = 1991.9615
$ws.Cells.Item(132, 10).Value = 2786.625
$ws.Cells.Item(132, 11).Value = 5975.8845
$ws.Cells.Item(132, 12).Value = 8359.875
$ws.Cells.Item(132, 13).Value = -3445.8845
$ws.Cells.Item(132, 14).Value = -13419.875
$ws.Cells.Item(136, 8).Value = 1569161.9
$ws.Cells.Item(136, 9).Value = 2273883
$ws.Cells.Item(136, 10).Value = 3114.5
$ws.Cells.Item(136, 11).Value = 6821649
$ws.Cells.Item(136, 12).Value = 9343.5
$ws.Cells.Item(136, 13).Value = -6819099
$ws.Cells.Item(136, 14).Value = -14443.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 1512.5
$ws.Cells.Item(17, 9).Value = 200
$ws.Cells.Item(17, 10).Value = 1700
$ws.Cells.Item(17, 11).Value = 600
$ws.Cells.Item(17, 12).Value = 5100
$ws.Cells.Item(17, 13).Value = -431
$ws.Cells.Item(17, 14).Value = -5438
$ws.Cells.Item(92, 8).Value = 426.2
$ws.Cells.Item(92, 9).Value = 596
$ws.Cells.Item(92, 10).Value = 313
$ws.Cells.Item(92, 11).Value = 1788
$ws.Cells.Item(92, 12).Value = 939
$ws.Cells.Item(92, 13).Value = -540
$ws.Cells.Item(92, 14).Value = -3435
$ws.Cells.Item(97, 8).Value = 13433.777
$ws.Cells.Item(97, 10).Value = 19509
$ws.Cells.Item(97, 12).Value = 58527
$ws.Cells.Item(97, 14).Value = -59519
$ws.Cells.Item(98, 8).Value = 400
$ws.Cells.Item(98, 9).Value = 400
$ws.Cells.Item(98, 11).Value = 1200
$ws.Cells.Item(98, 13).Value = 298

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5845.775
$ws.Cells.Item(70, 9).Value = 5194.25
$ws.Cells.Item(70, 10).Value = 6125
$ws.Cells.Item(70, 11).Value = 5194.25
$ws.Cells.Item(70, 12).Value = 6125
$ws.Cells.Item(70, 13).Value = -4924.25
$ws.Cells.Item(70, 14).Value = -6665
$ws.Cells.Item(73, 8).Value = 5845.775
$ws.Cells.Item(73, 9).Value = 5194.25
$ws.Cells.Item(73, 10).Value = 6125
$ws.Cells.Item(73, 11).Value = 5194.25
$ws.Cells.Item(73, 12).Value = 6125
$ws.Cells.Item(73, 13).Value = -4258.25
$ws.Cells.Item(73, 14).Value = -7997
$ws.Cells.Item(132, 8).Value = 1852.8431
$ws.Cells.Item(132, 9).Value = 1693.0435
$ws.Cells.Item(132, 10).Value = 3323
$ws.Cells.Item(132, 11).Value = 5079.1305
$ws.Cells.Item(132, 12).Value = 9969
$ws.Cells.Item(132, 13).Value = -2549.1305
$ws.Cells.Item(132, 14).Value = -15029
$ws.Cells.Item(135, 8).Value = 53973.168
$ws.Cells.Item(135, 10).Value = 53973.168
$ws.Cells.Item(135, 12).Value = 53973.168
$ws.Cells.Item(135, 14).Value = -64113.168

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 262.5238
$ws.Cells.Item(93, 9).Value = 232.15384
$ws.Cells.Item(93, 11).Value = 232.15384
$ws.Cells.Item(93, 13).Value = 1015.84616
$ws.Cells.Item(132, 8).Value = 8884.532999999999
$ws.Cells.Item(132, 9).Value = 13098.9375
$ws.Cells.Item(132, 10).Value = 4068.0715
$ws.Cells.Item(132, 11).Value = 39296.8125
$ws.Cells.Item(132, 12).Value = 12204.2145
$ws.Cells.Item(132, 13).Value = -36766.8125
$ws.Cells.Item(132, 14).Value = -17264.2145
$ws.Cells.Item(136, 8).Value = 4297.3
$ws.Cells.Item(136, 9).Value = 2205.9092
$ws.Cells.Item(136, 11).Value = 6617.7276
$ws.Cells.Item(136, 13).Value = -4067.7276
$ws.Cells.Item(141, 8).Value = 59147.5
$ws.Cells.Item(141, 10).Value = 59147.5
$ws.Cells.Item(141, 12).Value = 59147.5
$ws.Cells.Item(141, 14).Value = -69507.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 10008750
$ws.Cells.Item(5, 10).Value = 13341667
$ws.Cells.Item(5, 12).Value = 13341667
$ws.Cells.Item(5, 14).Value = -13341891
$ws.Cells.Item(132, 8).Value = 1712.9318
$ws.Cells.Item(132, 9).Value = 971.1429000000001
$ws.Cells.Item(132, 10).Value = 3011.0625
$ws.Cells.Item(132, 11).Value = 2913.4287
$ws.Cells.Item(132, 12).Value = 9033.1875
$ws.Cells.Item(132, 13).Value = -383.4287000000004
$ws.Cells.Item(132, 14).Value = -14093.1875
